$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text is numeric-looking need NumberFormat "@" forced
# before assignment (otherwise Excel auto-coerces to a number), then
# ClearFormats() afterwards so the cell keeps style index 0 (unstyled),
# matching the source file exactly while the value stays a text string.

$ws.Range('D2').Value = '26.462.97'
$ws.Range('E2').Value = '  -0.87%  '
$ws.Range('D3').Value = '1.627.10'
$ws.Range('E3').Value = '  -0.63%  '
$ws.Range('E4').Value = '  +0.19%  '
$ws.Range('E5').Value = '  -0.16%  '
$ws.Range('E6').Value = '  +1.39%  '
$ws.Range('E7').Value = '  +0.22%  '
$ws.Range('E8').Value = '  +0.47%  '
$ws.Range('E9').Value = '  -1.30%  '
$ws.Range('E10').Value = '  -0.99%  '
$ws.Range('E11').Value = '  +1.22%  '
$ws.Range('D12').Value = '1.853.68'
$ws.Range('E12').Value = '  -0.62%  '
$ws.Range('D13').Value = '1.643.89'
$ws.Range('E13').Value = '  +0.39%  '
$ws.Range('E14').Value = '  +1.73%  '
$ws.Range('E15').Value = '  -0.21%  '
$ws.Range('E16').Value = '  +2.86%  '
$ws.Range('D17').Value = '26.461.42'
$ws.Range('E17').Value = '  -0.83%  '
$ws.Range('E18').Value = '  +0.06%  '
$ws.Range('E19').Value = '  +2.17%  '
$ws.Range('E20').Value = '  +0.21%  '
$ws.Range('E21').Value = '  -0.42%  '
$ws.Range('E22').Value = '  +2.16%  '
$ws.Range('E23').Value = '  -1.22%  '
$ws.Range('E24').Value = '  +2.54%  '
$ws.Range('E25').Value = '  +1.28%  '
$ws.Range('E26').Value = '  +0.26%  '
$ws.Range('E27').Value = '  -0.23%  '
$ws.Range('E28').Value = '  +2.04%  '
$ws.Range('E29').Value = '  +0.84%  '
$ws.Range('E30').Value = '  -1.31%  '
$ws.Range('E31').Value = '  -1.15%  '
$ws.Range('E32').Value = '  +2.82%  '
$ws.Range('E33').Value = '  -0.61%  '
$ws.Range('E34').Value = '  -1.04%  '
$ws.Range('E35').Value = '  -1.07%  '
$ws.Range('D36').Value = '1.217.56'
$ws.Range('E36').Value = '  +4.63%  '
$ws.Range('E37').Value = '  +3.50%  '
$ws.Range('E38').Value = '  +0.25%  '
$ws.Range('E39').Value = '  -1.66%  '
$ws.Range('E40').Value = '  +0.45%  '
$ws.Range('E41').Value = '  -3.14%  '
$ws.Range('E42').Value = '  -0.42%  '
$ws.Range('E43').Value = '  +0.09%  '
$ws.Range('D44').Value = '1.763.86'
$ws.Range('E44').Value = '  -0.58%  '
$ws.Range('E45').Value = '  +0.61%  '
$ws.Range('E46').Value = '  +1.83%  '
$ws.Range('E47').Value = '  +0.42%  '
$ws.Range('B48').Value = 'Cronos'
$ws.Range('C48').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('E48').Value = '  -0.62%  '
$ws.Range('B49').Value = 'EnergySwap'
$ws.Range('C49').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('E49').Value = '  -0.10%  '
$ws.Range('B50').Value = 'Mantle'
$ws.Range('C50').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('E50').Value = '  -0.86%  '
$ws.Range('B51').Value = 'USDD'
$ws.Range('C51').Value = 'https://coinranking.com/coin/z2PZIKQL7+usdd-usdd'
$ws.Range('E51').Value = '  +0.32%  '

# Numeric-looking text values in column D (force text, then strip the
# temporary number-format style so the cell ends up unstyled again).
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '213.17'
$ws.Range('D5').ClearFormats()
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.0623'
$ws.Range('D8').ClearFormats()
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '18.94'
$ws.Range('D10').ClearFormats()
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0845'
$ws.Range('D11').ClearFormats()
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '64.84'
$ws.Range('D16').ClearFormats()
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '215.08'
$ws.Range('D19').ClearFormats()
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '4.30'
$ws.Range('D21').ClearFormats()
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.26'
$ws.Range('D22').ClearFormats()
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '148.05'
$ws.Range('D25').ClearFormats()
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '6.83'
$ws.Range('D28').ClearFormats()
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.39'
$ws.Range('D35').ClearFormats()
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.795'
$ws.Range('D39').ClearFormats()
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.504'
$ws.Range('D40').ClearFormats()
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.792'
$ws.Range('D42').ClearFormats()
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '92.96'
$ws.Range('D45').ClearFormats()
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '54.84'
$ws.Range('D47').ClearFormats()
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0509'
$ws.Range('D48').ClearFormats()
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '7.51'
$ws.Range('D49').ClearFormats()
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.406'
$ws.Range('D50').ClearFormats()
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.00'
$ws.Range('D51').ClearFormats()
